# Atualização Fases Liberta e Sula
# For each "Jogo" worksheet, the "Total Cartola" (G) values for the two
# teams get swapped with the opponent's "Cartola Sofrido" (H) value, and
# "Saldo Cartola" (I) is recomputed as Total Cartola - Cartola Sofrido.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $h2 = $ws.Range("H2").Value2
    $h3 = $ws.Range("H3").Value2

    # New "Total Cartola" is the opposing team's "Cartola Sofrido"
    $newG2 = $h3
    $newG3 = $h2

    $ws.Range("G2").Value = $newG2
    $ws.Range("G3").Value = $newG3

    # Recompute "Saldo Cartola" = Total Cartola - Cartola Sofrido
    $ws.Range("I2").Value = $newG2 - $h2
    $ws.Range("I3").Value = $newG3 - $h3
}
